$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the cryptos table (coin name / link / price / 1h volume %) cell by cell.
# Price cells (column D) are prefixed with a leading apostrophe so values like
# "309.13" are stored as literal text (matching the source data) instead of being
# auto-converted to a number by Excel's normal text-to-number inference.

# Row 2
$ws.Range("D2").Value = "'44.584.98"
$ws.Range("E2").Value = "  +1.53%  "

# Row 3
$ws.Range("D3").Value = "'2.250.71"
$ws.Range("E3").Value = "  +0.88%  "

# Row 4
$ws.Range("E4").Value = "  +0.08%  "

# Row 5
$ws.Range("D5").Value = "'309.13"
$ws.Range("E5").Value = "  +2.05%  "

# Row 6
$ws.Range("D6").Value = "'95.71"
$ws.Range("E6").Value = "  +1.46%  "

# Row 7
$ws.Range("D7").Value = "'0.574"
$ws.Range("E7").Value = "  +1.50%  "

# Row 8
$ws.Range("E8").Value = "  +0.18%  "

# Row 9
$ws.Range("E9").Value = "  +2.29%  "

# Row 10
$ws.Range("D10").Value = "'35.25"
$ws.Range("E10").Value = "  +2.84%  "

# Row 11
$ws.Range("D11").Value = "'0.0810"
$ws.Range("E11").Value = "  +1.09%  "

# Row 12
$ws.Range("D12").Value = "'7.31"
$ws.Range("E12").Value = "  +2.83%  "

# Row 13
$ws.Range("D13").Value = "'0.104"
$ws.Range("E13").Value = "  +1.37%  "

# Row 14
$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").Value = "'0.843"
$ws.Range("E14").Value = "  +4.30%  "

# Row 15
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "'2.254.29"
$ws.Range("E15").Value = "  -0.36%  "

# Row 16
$ws.Range("D16").Value = "'13.68"
$ws.Range("E16").Value = "  +2.99%  "

# Row 17
$ws.Range("D17").Value = "'44.269.83"
$ws.Range("E17").Value = "  +1.19%  "

# Row 18
$ws.Range("D18").Value = "'0.0₃0967"
$ws.Range("E18").Value = "  +1.93%  "

# Row 19
$ws.Range("D19").Value = "'6.42"
$ws.Range("E19").Value = "  +4.97%  "

# Row 20
$ws.Range("D20").Value = "'12.19"
$ws.Range("E20").Value = "  +1.72%  "

# Row 21
$ws.Range("D21").Value = "'65.96"
$ws.Range("E21").Value = "  +2.55%  "

# Row 22
$ws.Range("D22").Value = "'240.84"
$ws.Range("E22").Value = "  +2.25%  "

# Row 23
$ws.Range("E23").Value = "  +4.08%  "

# Row 24
$ws.Range("E24").Value = "  +4.23%  "

# Row 25
$ws.Range("E25").Value = "  +0.26%  "

# Row 26
$ws.Range("E26").Value = "  +5.87%  "

# Row 27
$ws.Range("D27").Value = "'9.89"
$ws.Range("E27").Value = "  +1.38%  "

# Row 28
$ws.Range("D28").Value = "'37.76"
$ws.Range("E28").Value = "  +5.11%  "

# Row 29
$ws.Range("E29").Value = "  +3.17%  "

# Row 30
$ws.Range("D30").Value = "'20.18"
$ws.Range("E30").Value = "  +1.39%  "

# Row 31
$ws.Range("D31").Value = "'152.74"
$ws.Range("E31").Value = "  +0.50%  "

# Row 32
$ws.Range("D32").Value = "'0.0803"
$ws.Range("E32").Value = "  +0.59%  "

# Row 33
$ws.Range("E33").Value = "  +0.10%  "

# Row 34
$ws.Range("E34").Value = "  -1.19%  "

# Row 35
$ws.Range("D35").Value = "'0.110"
$ws.Range("E35").Value = "  +2.20%  "

# Row 36
$ws.Range("E36").Value = "  +2.56%  "

# Row 37
$ws.Range("D37").Value = "'1.81"
$ws.Range("E37").Value = "  +3.16%  "

# Row 38
$ws.Range("D38").Value = "'3.46"
$ws.Range("E38").Value = "  +5.22%  "

# Row 39
$ws.Range("D39").Value = "'14.61"
$ws.Range("E39").Value = "  +0.30%  "

# Row 40
$ws.Range("E40").Value = "  +2.16%  "

# Row 41
$ws.Range("E41").Value = "  +3.29%  "

# Row 42
$ws.Range("E42").Value = "  +0.19%  "

# Row 43
$ws.Range("D43").Value = "'1.751.34"
$ws.Range("E43").Value = "  +2.02%  "

# Row 44
$ws.Range("D44").Value = "'81.54"
$ws.Range("E44").Value = "  -2.81%  "

# Row 45
$ws.Range("D45").Value = "'0.195"
$ws.Range("E45").Value = "  +6.23%  "

# Row 46
$ws.Range("B46").Value = "ordi"
$ws.Range("C46").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D46").Value = "'71.38"
$ws.Range("E46").Value = "  +5.34%  "

# Row 47
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "'100.11"
$ws.Range("E47").Value = "  +1.27%  "

# Row 48
$ws.Range("D48").Value = "'55.82"
$ws.Range("E48").Value = "  +4.61%  "

# Row 49
$ws.Range("B49").Value = "FraxShare"
$ws.Range("C49").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D49").Value = "'8.20"
$ws.Range("E49").Value = "  +2.92%  "

# Row 50
$ws.Range("B50").Value = "THORChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D50").Value = "'4.89"
$ws.Range("E50").Value = "  +1.42%  "

# Row 51
$ws.Range("E51").Value = "  +6.36%  "

# The apostrophe prefix above applies a "quote prefix" style to force text;
# reset column D (Price) back to the Normal style so the cell formatting
# matches the original (unstyled) Price cells.
$ws.Range("D2:D51").Style = "Normal"
